$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45185
}
